$d = $word.ActiveDocument

# Word's line-break character (maps to <w:br/> when assigned into Range.Text)
$br = [char]11

# Helper: replace the whole (non-paragraph-mark) text of a paragraph, preserving the
# run formatting of that paragraph's first run.
function Set-ParaText($paraIndex, [string]$newText) {
    $para = $d.Paragraphs($paraIndex)
    $rng = $para.Range
    $rng.MoveEnd(1, -1) # exclude the trailing paragraph mark
    $rng.Text = $newText
}

# --- Paragraph 1: Title ---
Set-ParaText 1 "Unveiling the Enigmatic Symphony of Politics: A High School Perspective"

# --- Paragraph 2: Author ---
Set-ParaText 2 "Prof. Julian Williamson"

# --- Paragraph 3: Email line ---
Set-ParaText 3 "politics.simplified@schoolconnect.org"

# --- Paragraph 5: Main body ---
$body = (
    "In the world of governance, politics stands as an enigmatic symphony, a tapestry woven with intrigue, power dynamics, and human interactions." +
    " At the heart of every political system lies the quest for order, stability, and the pursuit of a harmonious society." +
    " Yet, within this intricate symphony, there exist layers of complexity, challenges, and paradoxes that captivate the minds of students in high schools and beyond." +
    $br + $br +
    "Politics, in essence, is the art of negotiation, compromise, and decision-making." +
    " It delves into the study of how power is distributed, exercised, and contested within societies." +
    " Through the lens of politics, we seek to understand the motives of leaders, the influence of institutions, and the impact of policies on the lives of individuals and communities." +
    " It is a dynamic field where competing interests, ideologies, and values intertwine, shaping the course of nations and the lives of its citizens." +
    $br + $br +
    "The study of politics provides a crucial foundation for responsible citizenship." +
    " As future leaders, voters, and decision-makers, high school students play a pivotal role in shaping the political landscape." +
    " Understanding the nuances of politics equips them with the knowledge and skills necessary to navigate the complexities of governance, advocate for change, and make informed decisions that contribute to the betterment of society."
)
Set-ParaText 5 $body

# --- Paragraph 7: Summary body (paragraph 6 "Summary" heading is unchanged) ---
$summary = (
    "Politics, a multifaceted and dynamic field of study, presents a symphony of intrigue, power dynamics, and human interactions." +
    " It delves into the art of negotiation, compromise, and decision-making, seeking to understand the distribution and exercise of power within societies." +
    " Politics provides a foundation for responsible citizenship, empowering high school students to navigate the complexities of governance, advocate for change, and contribute to the betterment of society."
)
Set-ParaText 7 $summary

# --- Add a new empty paragraph at the very end of the document (after the Summary body) ---
$endPos = $d.Content.End
$endRng = $d.Range($endPos, $endPos)
$endRng.InsertParagraphAfter()
